$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 20.90688753949975
$ws.Range("C2").Value = 3.997979950226144
$ws.Range("E2").Value = 10.14393310931719
$ws.Range("F2").Value = 16.86991607391233
$ws.Range("G2").Value = 3.764652333072513
$ws.Range("K2").Value = 17.07269501980243
$ws.Range("L2").Value = 10.46822765296803
$ws.Range("N2").Value = 25.03386987323274
$ws.Range("B3").Value = 20.75441429926147
$ws.Range("C3").Value = 3.80623211594624
$ws.Range("E3").Value = 10.15001766818052
$ws.Range("F3").Value = 15.89584955866808
$ws.Range("G3").Value = 3.767945280939767
$ws.Range("K3").Value = 16.97529159512723
$ws.Range("L3").Value = 10.46364078319557
$ws.Range("N3").Value = 25.05826834335842
$ws.Range("B4").Value = 20.66570705687212
$ws.Range("C4").Value = 3.682616775326101
$ws.Range("E4").Value = 10.15539604784312
$ws.Range("F4").Value = 15.26997757108491
$ws.Range("G4").Value = 3.770071436352936
$ws.Range("K4").Value = 16.91947593065195
$ws.Range("L4").Value = 10.46277526412736
$ws.Range("N4").Value = 25.07497701631549
$ws.Range("B5").Value = 20.63082581103081
$ws.Range("C5").Value = 3.630791559048212
$ws.Range("E5").Value = 10.15800101487925
$ws.Range("F5").Value = 15.00819731993403
$ws.Range("G5").Value = 3.770964178342714
$ws.Range("K5").Value = 16.89775209215
$ws.Range("L5").Value = 10.46291384799606
$ws.Range("N5").Value = 25.08221998782327
$ws.Range("B6").Value = 20.62511127539577
$ws.Range("C6").Value = 3.622099307181902
$ws.Range("E6").Value = 10.15845853047861
$ws.Range("F6").Value = 14.96433081551593
$ws.Range("G6").Value = 3.771114009744205
$ws.Range("K6").Value = 16.89420705716173
$ws.Range("L6").Value = 10.46296654909204
$ws.Range("N6").Value = 25.08344888015277
$ws.Range("B7").Value = 20.66523146225206
$ws.Range("C7").Value = 3.681923676416929
$ws.Range("E7").Value = 10.155429505946
$ws.Range("F7").Value = 15.26647399323137
$ws.Range("G7").Value = 3.770083369502166
$ws.Range("K7").Value = 16.91917879675009
$ws.Range("L7").Value = 10.46277514311865
$ws.Range("N7").Value = 25.07507294077299
$ws.Range("B8").Value = 20.85331551164203
$ws.Range("C8").Value = 3.933101192837127
$ws.Range("E8").Value = 10.14569034617835
$ws.Range("F8").Value = 16.53996406344768
$ws.Range("G8").Value = 3.765766157927485
$ws.Range("K8").Value = 17.03829343608913
$ws.Range("L8").Value = 10.46624171559004
$ws.Range("N8").Value = 25.04192356970257
$ws.Range("B9").Value = 21.25958528409797
$ws.Range("C9").Value = 4.378010401369797
$ws.Range("E9").Value = 10.13961045426548
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 3.758123073558493
$ws.Range("K9").Value = 17.30267086466359
$ws.Range("L9").Value = 10.48847727322421
$ws.Range("N9").Value = 24.99064863533861
$ws.Range("B10").Value = 21.57869935840002
$ws.Range("C10").Value = 4.674812185435102
$ws.Range("E10").Value = 10.14305702170447
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 3.753003233820488
$ws.Range("K10").Value = 17.51441960340582
$ws.Range("L10").Value = 10.51415607236994
$ws.Range("N10").Value = 24.96137699703131
$ws.Range("B11").Value = 21.72785009081207
$ws.Range("C11").Value = 4.803156534719953
$ws.Range("E11").Value = 10.14633598241523
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 3.750780371392715
$ws.Range("K11").Value = 17.61425593503502
$ws.Range("L11").Value = 10.52784671203778
$ws.Range("N11").Value = 24.9498908449589
$ws.Range("B12").Value = 21.78485849471003
$ws.Range("C12").Value = 4.850788678130723
$ws.Range("E12").Value = 10.14782282809634
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 3.749953798856196
$ws.Range("K12").Value = 17.65253962174751
$ws.Range("L12").Value = 10.53331776675391
$ws.Range("N12").Value = 24.94580493691137
$ws.Range("B13").Value = 21.7725579405893
$ws.Range("C13").Value = 4.840573473345378
$ws.Range("E13").Value = 10.14749172047605
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 3.750131142475366
$ws.Range("K13").Value = 17.644273732483
$ws.Range("L13").Value = 10.53212676396936
$ws.Range("N13").Value = 24.94667317587246
$ws.Range("B14").Value = 21.73252987827893
$ws.Range("C14").Value = 4.807094728737246
$ws.Range("E14").Value = 10.14645339794511
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 3.750712065118237
$ws.Range("K14").Value = 17.61739613834113
$ws.Range("L14").Value = 10.52829108660568
$ws.Range("N14").Value = 24.94954940784966
$ws.Range("B15").Value = 21.70807898889448
$ws.Range("C15").Value = 4.786461581427422
$ws.Range("E15").Value = 10.14584929602898
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 3.751069870837489
$ws.Range("K15").Value = 17.60099424708483
$ws.Range("L15").Value = 10.52597889085025
$ws.Range("N15").Value = 24.95134553381692
$ws.Range("B16").Value = 21.56902838196234
$ws.Range("C16").Value = 4.66628944371977
$ws.Range("E16").Value = 10.14287709340002
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 3.753150631890674
$ws.Range("K16").Value = 17.50796342216508
$ws.Range("L16").Value = 10.51330160404864
$ws.Range("N16").Value = 24.96216450837301
$ws.Range("B17").Value = 21.48471430107307
$ws.Range("C17").Value = 4.590851229239498
$ws.Range("E17").Value = 10.14149141261487
$ws.Range("F17").Value = 20.20408069617459
$ws.Range("G17").Value = 3.754454241263337
$ws.Range("K17").Value = 17.45177229752534
$ws.Range("L17").Value = 10.50603761563322
$ws.Range("N17").Value = 24.96927060996213
$ws.Range("B18").Value = 21.43659681830347
$ws.Range("C18").Value = 4.54683355458184
$ws.Range("E18").Value = 10.14085552343346
$ws.Range("F18").Value = 19.95656407809808
$ws.Range("G18").Value = 3.755214042324805
$ws.Range("K18").Value = 17.41978477938743
$ws.Range("L18").Value = 10.50204882891392
$ws.Range("N18").Value = 24.9735300361997
$ws.Range("B19").Value = 21.42037128902776
$ws.Range("C19").Value = 4.531822464245597
$ws.Range("E19").Value = 10.14066791906735
$ws.Range("F19").Value = 19.87204792380562
$ws.Range("G19").Value = 3.755473018101692
$ws.Range("K19").Value = 17.40901220677522
$ws.Range("L19").Value = 10.50073086435645
$ws.Range("N19").Value = 24.97500176080593
$ws.Range("B20").Value = 21.49365090429801
$ws.Range("C20").Value = 4.598946787765819
$ws.Range("E20").Value = 10.14162225158957
$ws.Range("F20").Value = 20.24955283636157
$ws.Range("G20").Value = 3.754314435581497
$ws.Range("K20").Value = 17.45771975725318
$ws.Range("L20").Value = 10.50679130640704
$ws.Range("N20").Value = 24.96849632902349
$ws.Range("B21").Value = 21.74427311385483
$ws.Range("C21").Value = 4.816954613636132
$ws.Range("E21").Value = 10.14675173213681
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 3.750541022874912
$ws.Range("K21").Value = 17.62527799063812
$ws.Range("L21").Value = 10.52940995586065
$ws.Range("N21").Value = 24.94869742914835
$ws.Range("B22").Value = 21.91112704685464
$ws.Range("C22").Value = 4.95378384946221
$ws.Range("E22").Value = 10.15153270814067
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 3.74816330203755
$ws.Range("K22").Value = 17.7375590446742
$ws.Range("L22").Value = 10.5458624839241
$ws.Range("N22").Value = 24.93729469530083
$ws.Range("B23").Value = 21.82180913346091
$ws.Range("C23").Value = 4.881275239741512
$ws.Range("E23").Value = 10.14885062189226
$ws.Range("F23").Value = 21.82633154475857
$ws.Range("G23").Value = 3.749424275715606
$ws.Range("K23").Value = 17.67738803481143
$ws.Range("L23").Value = 10.53692946101503
$ws.Range("N23").Value = 24.94323973620724
$ws.Range("B24").Value = 21.48960955429667
$ws.Range("C24").Value = 4.595288801062564
$ws.Range("E24").Value = 10.14156259849868
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 3.754377609518395
$ws.Range("K24").Value = 17.45502991966177
$ws.Range("L24").Value = 10.50644997899591
$ws.Range("N24").Value = 24.96884583933344
$ws.Range("B25").Value = 21.14589021792937
$ws.Range("C25").Value = 4.262871585804485
$ws.Range("E25").Value = 10.13986338287135
$ws.Range("F25").Value = 18.34778573295691
$ws.Range("G25").Value = 3.760103265981166
$ws.Range("K25").Value = 17.22797793960634
$ws.Range("L25").Value = 10.48081533344863
$ws.Range("N25").Value = 25.00304681839252
